$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" (Changed) date column (C) for rows 2-5 from 45183 to 45184
$ws.Range("C2:C5").Value = 45184
